$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new standalone "ListView" locator string first (row 18) so that it is
# appended to the shared-strings table before the combined Abwesenheiten locator
# that reuses it as a prefix.
$ws.Range("E18").Value = "//android.webkit.WebView/android.widget.ListView"

# Update the Abwesenheiten mobile locator (D2) to the new, WebView/ListView-scoped xpath.
# (content-desc is prefixed with the same private-use icon glyph character as the
# original locator text, followed by a space and the visible label.)
$icon = [char]0xEE72
$ws.Range("D2").Value = "//android.webkit.WebView/android.widget.ListView//android.view.View[@content-desc=`"$icon Abwesenheiten`"]/android.widget.TextView[@text=`"Abwesenheiten`"]"

# Widen column D to fit the new, longer locator text.
$ws.Columns.Item(4).ColumnWidth = 109.67

# Update the active selection/view to D2 (also clears the old topLeftCell scroll position).
$ws.Range("D2").Select()
